$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 ("Separate roll sounds from crit sounds"): record a Completed Version of 1.7.2
$ws.Range("C32").Value = "1.7.2"

# Row 35 ("Better icons for saved rolls"): record a Completed Version of 1.7.2
$ws.Range("C35").Value = "1.7.2"

# Row 38: rename the "Allow for different dice images" feature to "Override dice icon"
# (description and requestor stay the same)
$ws.Range("A38").Value = "Override dice icon"

# New row 39: a follow-on feature request from the same requestor (Rae)
$ws.Range("A39").Value = "Dice icon themes"
$ws.Range("B39").Value = "I want to have all of my dice fit a set theme and change as a set."
$ws.Range("D39").Value = "Rae - urbanchika@gmail.com"

# Update the saved selection to match what was left selected after the edit
$ws.Range("C42").Select()
